# Auto-generated edit script: updates currentAveragePrice / Leve profit
# columns across all 8 profession sheets per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17
$ws.Range("H17").Value = 149525.23
$ws.Range("J17").Value = 149525.23
$ws.Range("L17").Value = 448575.6900000001
$ws.Range("N17").Value = -448911.6900000001

# ALC row 33
$ws.Range("H33").Value = 217
$ws.Range("I33").Value = 206
$ws.Range("J33").Value = 231.3
$ws.Range("K33").Value = 206
$ws.Range("L33").Value = 231.3
$ws.Range("M33").Value = 23
$ws.Range("N33").Value = -689.3

# ALC row 40
$ws.Range("H40").Value = 4356.4287
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 4749.1665
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 4749.1665
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -5099.1665

# ALC row 69
$ws.Range("H69").Value = 3150
$ws.Range("I69").Value = 1900
$ws.Range("K69").Value = 5700
$ws.Range("M69").Value = -4826

# ALC row 72
$ws.Range("H72").Value = 3150
$ws.Range("I72").Value = 1900
$ws.Range("K72").Value = 17100
$ws.Range("M72").Value = -12732

# ALC row 113
$ws.Range("H113").Value = 3339
$ws.Range("J113").Value = 4221.5557
$ws.Range("L113").Value = 4221.5557
$ws.Range("N113").Value = -10729.5557

# ALC row 132
$ws.Range("H132").Value = 27036284
$ws.Range("I132").Value = 33336376
$ws.Range("K132").Value = 100009128
$ws.Range("M132").Value = -100006598

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Range("H61").Value = 6489.6587
$ws.Range("I61").Value = 4140.778
$ws.Range("K61").Value = 4140.778
$ws.Range("M61").Value = -3928.778

# ARM row 74
$ws.Range("H74").Value = 2889.5898
$ws.Range("I74").Value = 2357.2285
$ws.Range("K74").Value = 2357.2285
$ws.Range("M74").Value = -1483.2285

# ARM row 77
$ws.Range("H77").Value = 2889.5898
$ws.Range("I77").Value = 2357.2285
$ws.Range("K77").Value = 11786.1425
$ws.Range("M77").Value = -7418.142500000002

# ARM row 132
$ws.Range("H132").Value = 2752
$ws.Range("I132").Value = 2609
$ws.Range("J132").Value = 4325
$ws.Range("K132").Value = 7827
$ws.Range("L132").Value = 12975
$ws.Range("M132").Value = -5297
$ws.Range("N132").Value = -18035

# ARM row 136
$ws.Range("H136").Value = 6489.6587
$ws.Range("I136").Value = 4140.778
$ws.Range("K136").Value = 12422.334
$ws.Range("M136").Value = -9872.334000000001

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 2451.0356
$ws.Range("I20").Value = 2231.6875
$ws.Range("K20").Value = 2231.6875
$ws.Range("M20").Value = -1984.6875

# BSM row 33
$ws.Range("H33").Value = 13498.75
$ws.Range("J33").Value = 13498.75
$ws.Range("L33").Value = 13498.75
$ws.Range("N33").Value = -14170.75

# BSM row 39
$ws.Range("H39").Value = 12474.25
$ws.Range("J39").Value = 12474.25
$ws.Range("L39").Value = 12474.25
$ws.Range("N39").Value = -13252.25

# BSM row 94
$ws.Range("H94").Value = 909.6875
$ws.Range("I94").Value = 462.16666
$ws.Range("K94").Value = 462.16666
$ws.Range("M94").Value = -11.16665999999998

# BSM row 99
$ws.Range("H99").Value = 2976.6538
$ws.Range("I99").Value = 2418.375
$ws.Range("J99").Value = 3869.9
$ws.Range("K99").Value = 2418.375
$ws.Range("L99").Value = 3869.9
$ws.Range("M99").Value = -920.375
$ws.Range("N99").Value = -6865.9

# BSM row 134
$ws.Range("H134").Value = 8783.902
$ws.Range("I134").Value = 3177.4285
$ws.Range("J134").Value = 9938.177
$ws.Range("K134").Value = 9532.2855
$ws.Range("L134").Value = 29814.531
$ws.Range("M134").Value = -6997.2855
$ws.Range("N134").Value = -34884.531

$ws = $wb.Worksheets.Item("CRP")
# CRP row 58
$ws.Range("H58").Value = 1654.1904
$ws.Range("I58").Value = 1547
$ws.Range("J58").Value = 1997.2
$ws.Range("K58").Value = 1547
$ws.Range("L58").Value = 1997.2
$ws.Range("M58").Value = -1344
$ws.Range("N58").Value = -2403.2

# CRP row 94
$ws.Range("H94").Value = 1650.6111
$ws.Range("I94").Value = 1129.8
$ws.Range("J94").Value = 1850.9231
$ws.Range("K94").Value = 1129.8
$ws.Range("L94").Value = 1850.9231
$ws.Range("M94").Value = -678.8
$ws.Range("N94").Value = -2752.9231

# CRP row 132
$ws.Range("H132").Value = 2108164.8
$ws.Range("I132").Value = 2225146.2
$ws.Range("K132").Value = 6675438.600000001
$ws.Range("M132").Value = -6672908.600000001

# CRP row 134
$ws.Range("H134").Value = 3162.4182
$ws.Range("I134").Value = 1287.7368
$ws.Range("K134").Value = 3863.2104
$ws.Range("M134").Value = -1328.2104

# CRP row 136
$ws.Range("H136").Value = 1654.1904
$ws.Range("I136").Value = 1547
$ws.Range("J136").Value = 1997.2
$ws.Range("K136").Value = 4641
$ws.Range("L136").Value = 5991.6
$ws.Range("M136").Value = -2091
$ws.Range("N136").Value = -11091.6

$ws = $wb.Worksheets.Item("CUL")
# CUL row 29
$ws.Range("H29").Value = 2303.9
$ws.Range("I29").Value = 572.25
$ws.Range("K29").Value = 1716.75
$ws.Range("M29").Value = -1439.75

$ws = $wb.Worksheets.Item("GSM")
# GSM row 3
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("M3").Value = -9884

# GSM row 132
$ws.Range("H132").Value = 23814024
$ws.Range("I132").Value = 37040948
$ws.Range("J132").Value = 5559.6
$ws.Range("K132").Value = 111122844
$ws.Range("L132").Value = 16678.8
$ws.Range("M132").Value = -111120314
$ws.Range("N132").Value = -21738.8

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93
$ws.Range("H93").Value = 4444.8335
$ws.Range("I93").Value = 4804.75
$ws.Range("K93").Value = 4804.75
$ws.Range("M93").Value = -3556.75

# LTW row 132
$ws.Range("H132").Value = 3131.0889
$ws.Range("I132").Value = 3121.5278
$ws.Range("K132").Value = 9364.5834
$ws.Range("M132").Value = -6834.5834

# LTW row 136
$ws.Range("H136").Value = 5309.45
$ws.Range("I136").Value = 4928.5
$ws.Range("J136").Value = 6198.3335
$ws.Range("K136").Value = 14785.5
$ws.Range("L136").Value = 18595.0005
$ws.Range("M136").Value = -12235.5
$ws.Range("N136").Value = -23695.0005

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws.Range("H81").Value = 49508.434
$ws.Range("I81").Value = 79841.30499999999
$ws.Range("K81").Value = 159682.61
$ws.Range("M81").Value = -158621.61

# WVR row 84
$ws.Range("H84").Value = 49508.434
$ws.Range("I84").Value = 79841.30499999999
$ws.Range("K84").Value = 798413.0499999999
$ws.Range("M84").Value = -793109.0499999999

# WVR row 100
$ws.Range("H100").Value = 2993.923
$ws.Range("J100").Value = 3909.6
$ws.Range("L100").Value = 7819.2
$ws.Range("N100").Value = -8901.200000000001

# WVR row 122
$ws.Range("H122").Value = 1861.7391
$ws.Range("I122").Value = 1570.55
$ws.Range("K122").Value = 4711.65
$ws.Range("M122").Value = -2261.65

# WVR row 136
$ws.Range("H136").Value = 9134.959999999999
$ws.Range("I136").Value = 13710.9375
$ws.Range("K136").Value = 41132.8125
$ws.Range("M136").Value = -38582.8125
